$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03833417482769
$ws.Range("D2").Value = 1.041792446057049
$ws.Range("E2").Value = 1.052032207462182
$ws.Range("F2").Value = 1.059771666515263
$ws.Range("I2").Value = 1.040317102431977
$ws.Range("J2").Value = 1.04343217696631
$ws.Range("K2").Value = 1.044570958553282
$ws.Range("L2").Value = 1.054782063451087
$ws.Range("M2").Value = 1.062500261231362
$ws.Range("N2").Value = 1.018396581411455
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039177860110964
$ws.Range("D3").Value = 1.042429439714194
$ws.Range("E3").Value = 1.052867238487455
$ws.Range("F3").Value = 1.06069112053527
$ws.Range("I3").Value = 1.040515125507956
$ws.Range("J3").Value = 1.043921142922674
$ws.Range("K3").Value = 1.045019108057099
$ws.Range("L3").Value = 1.055429802046257
$ws.Range("M3").Value = 1.063233744205571
$ws.Range("N3").Value = 1.018560634815528
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039724423549446
$ws.Range("D4").Value = 1.042842116663374
$ws.Range("E4").Value = 1.053408577159608
$ws.Range("F4").Value = 1.06128724151901
$ws.Range("I4").Value = 1.040642410101705
$ws.Range("J4").Value = 1.044237511790226
$ws.Range("K4").Value = 1.045308900607756
$ws.Range("L4").Value = 1.055849300118856
$ws.Range("M4").Value = 1.063708895597561
$ws.Range("N4").Value = 1.018666733691259
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039954351114867
$ws.Range("D5").Value = 1.043015723931967
$ws.Range("E5").Value = 1.053636397811042
$ws.Range("F5").Value = 1.061538129719203
$ws.Range("I5").Value = 1.04069571647186
$ws.Range("J5").Value = 1.044370506130032
$ws.Range("K5").Value = 1.045430682665484
$ws.Range("L5").Value = 1.056025743717313
$ws.Range("M5").Value = 1.063908776669779
$ws.Range("N5").Value = 1.018711324123824
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039992965850424
$ws.Range("D6").Value = 1.04304488021954
$ws.Range("E6").Value = 1.053674664032723
$ws.Range("F6").Value = 1.061580271256203
$ws.Range("I6").Value = 1.04070465485923
$ws.Range("J6").Value = 1.04439283601196
$ws.Range("K6").Value = 1.045451127633395
$ws.Range("L6").Value = 1.056055374410114
$ws.Range("M6").Value = 1.063942345007889
$ws.Range("N6").Value = 1.018718810246823
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039727495256347
$ws.Range("D7").Value = 1.042844435951135
$ws.Range("E7").Value = 1.053411620362135
$ws.Range("F7").Value = 1.06129059280566
$ws.Range("I7").Value = 1.040643123186748
$ws.Range("J7").Value = 1.044239288896621
$ws.Range("K7").Value = 1.045310528050705
$ws.Range("L7").Value = 1.055851657428373
$ws.Range("M7").Value = 1.063711565919502
$ws.Range("N7").Value = 1.018667329564551
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03861916821972
$ws.Range("D8").Value = 1.042007616612634
$ws.Range("E8").Value = 1.052314198301539
$ws.Range("F8").Value = 1.060082156385295
$ws.Range("I8").Value = 1.040384200757453
$ws.Range("J8").Value = 1.043597429717688
$ws.Range("K8").Value = 1.044722451419792
$ws.Range("L8").Value = 1.055000892687352
$ws.Range("M8").Value = 1.062748032712904
$ws.Range("N8").Value = 1.018452035070502
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036671145115757
$ws.Range("D9").Value = 1.040536934080086
$ws.Range("E9").Value = 1.050388274047347
$ws.Range("F9").Value = 1.057961799466998
$ws.Range("I9").Value = 1.039921469852767
$ws.Range("J9").Value = 1.042466254727147
$ws.Range("K9").Value = 1.043684782480096
$ws.Range("L9").Value = 1.05350461891097
$ws.Range("M9").Value = 1.061054360953965
$ws.Range("N9").Value = 1.018072260131101
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035375911124437
$ws.Range("D10").Value = 1.039559202198092
$ws.Range("E10").Value = 1.049109724548928
$ws.Range("F10").Value = 1.056554429848424
$ws.Range("I10").Value = 1.039608668931062
$ws.Range("J10").Value = 1.041712117086902
$ws.Range("K10").Value = 1.042992132175661
$ws.Range("L10").Value = 1.052509128359177
$ws.Range("M10").Value = 1.059928162264758
$ws.Range("N10").Value = 1.017818835835243
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034815897787729
$ws.Range("D11").Value = 1.039136500724814
$ws.Range("E11").Value = 1.048557401715347
$ws.Range("F11").Value = 1.055946517100101
$ws.Range("I11").Value = 1.039472208122711
$ws.Range("J11").Value = 1.041385577183917
$ws.Range("K11").Value = 1.042692015562946
$ws.Range("L11").Value = 1.052078568738936
$ws.Range("M11").Value = 1.059441217215251
$ws.Range("N11").Value = 1.017709048834058
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034608010383187
$ws.Range("D12").Value = 1.038979591894655
$ws.Range("E12").Value = 1.048352441174316
$ws.Range("F12").Value = 1.055720936578202
$ws.Range("I12").Value = 1.039421368545196
$ws.Range("J12").Value = 1.04126428792901
$ws.Range("K12").Value = 1.04258051089416
$ws.Range("L12").Value = 1.051918715477338
$ws.Range("M12").Value = 1.059260451868323
$ws.Range("N12").Value = 1.017668261605361
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034652597188499
$ws.Range("D13").Value = 1.039013244758809
$ws.Range("E13").Value = 1.048396396973609
$ws.Range("F13").Value = 1.055769314139098
$ws.Range("I13").Value = 1.039432280672652
$ws.Range("J13").Value = 1.041290304766859
$ws.Range("K13").Value = 1.042604430271427
$ws.Range("L13").Value = 1.051953001092866
$ws.Range("M13").Value = 1.05929922174842
$ws.Range("N13").Value = 1.017677010934103
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034798711160786
$ws.Range("D14").Value = 1.039123528515131
$ws.Range("E14").Value = 1.048540455588936
$ws.Range("F14").Value = 1.055927865936227
$ws.Range("I14").Value = 1.039468008803335
$ws.Range("J14").Value = 1.041375551324349
$ws.Range("K14").Value = 1.042682799116298
$ws.Range("L14").Value = 1.052065353664581
$ws.Range("M14").Value = 1.059426272888481
$ws.Range("N14").Value = 1.017705677498316
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03488875360835
$ws.Range("D15").Value = 1.039191491454446
$ws.Range("E15").Value = 1.04862924098206
$ws.Range("F15").Value = 1.056025584846366
$ws.Range("I15").Value = 1.039490001968288
$ws.Range("J15").Value = 1.041428074869967
$ws.Range("K15").Value = 1.042731081074847
$ws.Range("L15").Value = 1.052134587874302
$ws.Range("M15").Value = 1.05950456761153
$ws.Range("N15").Value = 1.017723338941314
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035413095017269
$ws.Range("D16").Value = 1.039587269616391
$ws.Range("E16").Value = 1.049146407888554
$ws.Range("F16").Value = 1.056594806514744
$ws.Range("I16").Value = 1.039617704038418
$ws.Range("J16").Value = 1.041733788707961
$ws.Range("K16").Value = 1.043012045933846
$ws.Range("L16").Value = 1.052537713725038
$ws.Range("M16").Value = 1.059960494223412
$ws.Range("N16").Value = 1.01782612096743
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035742224250808
$ws.Range("D17").Value = 1.039835709493505
$ws.Range("E17").Value = 1.049471161691346
$ws.Range("F17").Value = 1.056952263800569
$ws.Range("I17").Value = 1.039697536652926
$ws.Range("J17").Value = 1.041925557555999
$ws.Range("K17").Value = 1.043188236564931
$ws.Range("L17").Value = 1.052790717318991
$ws.Range("M17").Value = 1.060246675194112
$ws.Range("N17").Value = 1.017890579653807
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035934279890238
$ws.Range("D18").Value = 1.03998068411349
$ws.Range("E18").Value = 1.049660710162034
$ws.Range("F18").Value = 1.057160905993928
$ws.Range("I18").Value = 1.03974400359682
$ws.Range("J18").Value = 1.042037413643655
$ws.Range("K18").Value = 1.043290986734747
$ws.Range("L18").Value = 1.052938337633881
$ws.Range("M18").Value = 1.060413667764228
$ws.Range("N18").Value = 1.017928172208485
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035999779443604
$ws.Range("D19").Value = 1.040030127488727
$ws.Range("E19").Value = 1.049725362451439
$ws.Range("F19").Value = 1.057232071858078
$ws.Range("I19").Value = 1.039759830969027
$ws.Range("J19").Value = 1.042075553737882
$ws.Range("K19").Value = 1.043326018653931
$ws.Range("L19").Value = 1.05298868037501
$ws.Range("M19").Value = 1.060470619409421
$ws.Range("N19").Value = 1.017940989418129
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035706903499872
$ws.Range("D20").Value = 1.039809047629869
$ws.Range("E20").Value = 1.049436305731903
$ws.Range("F20").Value = 1.056913897170234
$ws.Range("I20").Value = 1.039688981507825
$ws.Range("J20").Value = 1.041904982506575
$ws.Range("K20").Value = 1.043169334918035
$ws.Range("L20").Value = 1.052763567493692
$ws.Range("M20").Value = 1.060215963633932
$ws.Range("N20").Value = 1.017883664371694
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034755680733672
$ws.Range("D21").Value = 1.039091049886916
$ws.Range("E21").Value = 1.048498028471431
$ws.Range("F21").Value = 1.055881170148432
$ws.Range("I21").Value = 1.039457491948669
$ws.Range("J21").Value = 1.041350448267377
$ws.Range("K21").Value = 1.04265972220781
$ws.Range("L21").Value = 1.052032266529503
$ws.Range("M21").Value = 1.059388856507824
$ws.Range("N21").Value = 1.017697236111598
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034158341793119
$ws.Range("D22").Value = 1.038640203255912
$ws.Range("E22").Value = 1.04790923555875
$ws.Range("F22").Value = 1.055233158345836
$ws.Range("I22").Value = 1.039311066433263
$ws.Range("J22").Value = 1.041001803930711
$ws.Range("K22").Value = 1.042339147097527
$ws.Range("L22").Value = 1.051572907599755
$ws.Range("M22").Value = 1.058869445488468
$ws.Range("N22").Value = 1.01757997839772
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034474932508172
$ws.Range("D23").Value = 1.038879149363761
$ws.Range("E23").Value = 1.04822125723009
$ws.Range("F23").Value = 1.0555765573785
$ws.Range("I23").Value = 1.039388772446855
$ws.Range("J23").Value = 1.041186625223185
$ws.Range("K23").Value = 1.042509104917924
$ws.Range("L23").Value = 1.051816380351854
$ws.Range("M23").Value = 1.059144735410762
$ws.Range("N23").Value = 1.017642142841643
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035722863179464
$ws.Range("D24").Value = 1.039821094781359
$ws.Range("E24").Value = 1.049452055253395
$ws.Range("F24").Value = 1.056931232956137
$ws.Range("I24").Value = 1.039692847512744
$ws.Range("J24").Value = 1.041914279484174
$ws.Range("K24").Value = 1.043177875817033
$ws.Range("L24").Value = 1.052775835184204
$ws.Range("M24").Value = 1.060229840655941
$ws.Range("N24").Value = 1.017886789105851
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037174154828436
$ws.Range("D25").Value = 1.040916668181653
$ws.Range("E25").Value = 1.050885227624768
$ws.Range("F25").Value = 1.058508877428793
$ws.Range("I25").Value = 1.040041860188395
$ws.Range("J25").Value = 1.042758699491904
$ws.Range("K25").Value = 1.043953203580112
$ws.Range("L25").Value = 1.053891090974873
$ws.Range("M25").Value = 1.061491708664644
$ws.Range("N25").Value = 1.01817048550731
